# Rev D of audio spectrum board
# Updates the BOM worksheet: quantity/part corrections, a couple of
# device/package and part-number swaps, expanded reference-designator
# lists for a few rows, and removal of the obsolete test-point line item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level corrections (row numbers as they exist BEFORE the
#     test-point row is deleted later in this script) -----------------

# Row 2: C1, C2 - 8pF caps switched from 0603 package to 0402 package
$ws.Range("A2").Value = 2
$ws.Range("C2").Value = "CAPACITOR0402_CAP"
$ws.Range("D2").Value = "0402_CAP"
$ws.Range("G2").Value = "490-8230-1-ND"

# Row 3: C12, C13 - qty correction
$ws.Range("A3").Value = 2

# Row 4: C14 - qty correction
$ws.Range("A4").Value = 1

# Row 5: C15-C19 - qty correction
$ws.Range("A5").Value = 5

# Row 6: C20-C23 - qty correction
$ws.Range("A6").Value = 4

# Row 7: C24 - qty correction
$ws.Range("A7").Value = 1

# Row 8: C25, C26, C27 -> also now used for C33, C34; qty correction
$ws.Range("A8").Value = 5
$ws.Range("E8").Value = "C25, C26, C27, C33, C34"

# Row 9: C28 - qty correction
$ws.Range("A9").Value = 1

# Row 11: C3, C4 - qty correction
$ws.Range("A11").Value = 2

# Row 12: C31 - qty correction
$ws.Range("A12").Value = 1

# Row 14: C5, C6 - qty correction
$ws.Range("A14").Value = 2

# Row 16: C9, C10, C11 - qty correction
$ws.Range("A16").Value = 3

# Row 17: FB1 ferrite bead switched from 0603 to 0402 package
$ws.Range("A17").Value = 1
$ws.Range("C17").Value = "FERRITE_BEAD0402"
$ws.Range("D17").Value = "0402_CAP"

# Row 18: J1 - qty correction
$ws.Range("A18").Value = 1

# Row 19: J2 - qty correction
$ws.Range("A19").Value = 1

# Row 21: L1 - qty correction
$ws.Range("A21").Value = 1

# Row 22: L2 - qty correction
$ws.Range("A22").Value = 1

# Row 23/24: LED reference designators renumbered (LED1..LED4 -> LED0..LED3)
$ws.Range("A23").Value = 2
$ws.Range("E23").Value = "LED0, LED2"
$ws.Range("G23").Value = "160-1447-1-ND"

$ws.Range("A24").Value = 2
$ws.Range("E24").Value = "LED1, LED3"
$ws.Range("G24").Value = "160-1446-1-ND"

# Row 25: P1, P2 - qty correction
$ws.Range("A25").Value = 2

# Row 26: Q1, Q2 - qty correction
$ws.Range("A26").Value = 2

# Row 27: R1 - qty correction
$ws.Range("A27").Value = 1

# Row 28: R10, R11, R13, R14 - qty correction
$ws.Range("A28").Value = 4

# Row 31: 20k resistors, add R32-R35 to the reference list; qty correction
$ws.Range("A31").Value = 8
$ws.Range("E31").Value = "R17, R18, R19, R20, R32, R33, R34, R35"

# Row 34: 100k resistors - qty correction
$ws.Range("A34").Value = 6

# Row 35: 200k resistors - qty correction
$ws.Range("A35").Value = 2

# Row 36: 1k resistors - qty correction
$ws.Range("A36").Value = 3

# Row 38: 750k resistor - qty correction
$ws.Range("A38").Value = 1

# Row 39: 1.5k resistors - qty correction
$ws.Range("A39").Value = 4

# Row 42: U1 flash memory - qty correction
$ws.Range("A42").Value = 1

# Row 43: U10 regulator - qty correction
$ws.Range("A43").Value = 1

# Row 44: U2 microphone - qty correction
$ws.Range("A44").Value = 1

# Row 45: U3, U4 spectrum analyzer - qty correction
$ws.Range("A45").Value = 2

# Row 46: U5-U8 op-amps - qty correction
$ws.Range("A46").Value = 4

# Row 47: U9 microcontroller - qty correction
$ws.Range("A47").Value = 1

# --- Remove the obsolete test-point line item (row 41), shifting all
#     rows below it up by one -----------------------------------------
$ws.Rows.Item(41).Delete()

# --- Cosmetic view state matching the refreshed sheet -----------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("J41").Select()
